# Apply "Added all symbiont density sample data" edit to the Daily sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Initials" column (D) to make room
# for the new "Time.Feeding" column. This shifts old D->E and old E->F,
# carrying all existing values/shared-strings along with it.
$ws.Columns("D:D").Insert()

# New column header and width (match neighboring "Visual.Inspection" column).
$ws.Range("D1").Value = "Time.Feeding"
$ws.Columns("D:D").ColumnWidth = $ws.Columns("C:C").ColumnWidth

# Backfill the new column with "NA" for all pre-existing daily rows.
$ws.Range("D2:D19").Value = "NA"

# Add the new row of data for 2022-10-17.
$ws.Range("A20").Value = 20221017
$ws.Range("B20").Value = "completed"
$ws.Range("C20").Value = "completed"
$ws.Range("F20").Value = "Started shade and feeding treatments"
$ws.Range("D20").Value = "13:53-14:55"
$ws.Range("E20").Value = "AH"

# Move the active selection to reflect the next empty row, as in the source file.
$ws.Range("A21").Select() | Out-Null
